$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing column A (url) shifts to column B.
$ws.Columns("A").Insert()

# Column widths: new A is narrower (company name), B keeps the original url width.
$ws.Columns("A").ColumnWidth = 35.33203125
$ws.Columns("B").ColumnWidth = 56.21875

# Header row
$ws.Range("A1").Value = "CompanyName"
$ws.Range("B1").Value = "URL"

# Data row
$ws.Range("A2").Value = "Trigent Software Ltd"
$ws.Range("B2").Value = "https://www.trigent.com/"

# Fix up the hyperlink: remove the old one (still anchored at A2) and add it back on B2.
$ws.Hyperlinks.Item(1).Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.trigent.com/", [type]::Missing, [type]::Missing, "https://www.trigent.com/")

# Header formatting: bold text on a themed fill.
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.Interior.ThemeColor = 3

# Selection / view bits
$ws.Range("A9").Select()
$excel.ActiveWindow.WindowState = -4143
